$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (SamplesTab): replace the "DISTINCT sample_id" query text with the new,
# shorter version that drops the Tumor / Analyte Type columns.
$newSamplesQuery = @"
SELECT
    DISTINCT (smp.sample_id) AS "Sample ID",
    sp.participant_id AS "Participant ID", 
    s.study_name AS "Study Name",
    s.phs_accession AS Accession
FROM 
    df_participant sp
JOIN 
    df_study s ON sp."study.phs_accession" = s.phs_accession
JOIN 
    df_sample smp ON smp."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_diagnosis d ON d."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_program p ON p.program_acronym = s."program.program_acronym"
JOIN
    df_file f1 ON f1."sample.sample_id" = smp.sample_id
JOIN
    df_genomic_info gi ON gi."file.file_id" = f1.file_id
WHERE 
   s.phs_accession = 'phs001437' AND f1.file_type = 'PDF'
ORDER BY 
    smp.sample_id ASC
LIMIT 100;
"@

$ws.Range("B3").Value2 = $newSamplesQuery

# Clear the trailing TSV/Web filename cells on rows 3 and 4 -- they now only
# appear on row 2.
$ws.Range("D3").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("D4").ClearContents()
$ws.Range("E4").ClearContents()

# Update the current selection to match the new state.
$ws.Range("B3").Select()
